# Updating some test models to comply with new xlsx and yml settings files
# structure.
#
# 1) Rename the shared header labels used across the "x", "a", "b" and
#    "products_data" sheets:
#       r_Names  -> resources_Name
#       p_Names  -> products_Name
#       pd_Names -> product_data_Name
#
# 2) Refresh the saved UI state (active sheet / selected range per sheet)
#    left over from the editing session.

$wb = $excel.ActiveWorkbook

# --- 1) Rename header labels (every cell that carries the old label) -------

$wsX  = $wb.Worksheets.Item("x")
$wsA  = $wb.Worksheets.Item("a")
$wsPD = $wb.Worksheets.Item("products_data")
$wsB  = $wb.Worksheets.Item("b")

# r_Names -> resources_Name
$wsX.Range("B1").Value = "resources_Name"
$wsA.Range("B1").Value = "resources_Name"
$wsB.Range("B1").Value = "resources_Name"

# p_Names -> products_Name
$wsX.Range("C1").Value = "products_Name"
$wsA.Range("C1").Value = "products_Name"
$wsPD.Range("B1").Value = "products_Name"

# pd_Names -> product_data_Name
$wsPD.Range("C1").Value = "product_data_Name"

# --- 2) Restore per-sheet selections, then leave "x" as the active tab -----

$wsA.Range("D2:D7").Select() | Out-Null
$wsPD.Range("D2:D7").Select() | Out-Null
$wsB.Range("C2:C4").Select() | Out-Null

# Activate "x" last so it becomes the saved active/selected sheet, matching
# the new saved selection of E17 on that sheet.
$wsX.Activate() | Out-Null
$wsX.Range("E17").Select() | Out-Null
